# Shift all timestamps in column A (rows 2-249) back by 3 hours (10800 seconds).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 249
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $v = $cell.Value()
    $cell.Value = $v - 10800
}
